# Weekly update: insert a new price record for the latest week, pushing the
# existing historical rows (previously rows 3-16) down by one row
# (becoming rows 4-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3 (row 2 - the most recent existing
# record - stays put; all rows from the old row 3 onward shift down by one).
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 'Macroferia Regional de Talca'
$ws.Range("C3").Value = 'Maule'
$ws.Range("D3").Value = '2022-05-10'
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 'Fruta'
$ws.Range("G3").Value = 100104
$ws.Range("H3").Value = 'Frutos de pepita'
$ws.Range("I3").Value = 100104001
$ws.Range("J3").Value = 'Granada'
$ws.Range("K3").Value = 'Wonderfull'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 17000
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 944
$ws.Range("T3").Value = 18
